$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-11, columns B-G (Q0..Q9 statistics table)
$data = @{
    2  = @(0.2143979625585434, 1.424134207528472, 7.525940274036349, 2.743344723879292, 2.796421391031461, 23)
    3  = @(0.4670168130507297, 1.853629581901463, 13.55051072646238, 3.681101835926627, 3.737282769615654, 22)
    4  = @(-0.8354076150457171, 1.386534768542286, 7.40374990259892, 2.720983260257019, 2.653513938212625, 21)
    5  = @(-0.04486636216424342, 0.5141898899062147, 0.5563440360261616, 0.7458847337398464, 0.7638758824313834, 20)
    6  = @(0.08415703371898316, 0.7459901986017395, 1.174992571369782, 1.08397074285692, 1.110312597020984, 19)
    7  = @(-0.2038590712083605, 0.6468456134312068, 0.8625708541187559, 0.9287469268421596, 0.9323664388718182, 18)
    8  = @(-0.1285815229469241, 0.5264431630486278, 0.5203200017941761, 0.7213321022900452, 0.7316238571562934, 17)
    9  = @(0.08851026789393576, 0.4260801985527543, 0.3054846069715224, 0.5527066192579227, 0.5634659783389385, 16)
    10 = @(0.007460085389865952, 0.3714499647231649, 0.2993920067773447, 0.5471672566750908, 0.5663192757087402, 15)
    11 = @(0.005873690370570317, 0.3581201956056038, 0.2817709375354674, 0.5308210032915686, 0.5524621921347255, 13)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
